# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "Datos actualizados a ..." timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 4 de Septiembre de 2020 a las 06:08"

# Row 6 - India
$ws.Range("B6").Value = 3936747
$ws.Range("C6").Value = 3623
$ws.Range("D6").Value = 3037151
$ws.Range("E6").Value = 831027

# Row 20 - Pakistan
$ws.Range("B20").Value = 297512
$ws.Range("C20").Value = 498
$ws.Range("D20").Value = 282268
$ws.Range("E20").Value = 8909
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = 6335

# Row 67 - Azerbaiyan
$ws.Range("D67").Value = 34264
$ws.Range("E67").Value = 2094

# Row 133 - Jamaica
$ws.Range("B133").Value = 2896
$ws.Range("C133").Value = 74
$ws.Range("D133").Value = 921
$ws.Range("E133").Value = 1946
$ws.Range("G133").Value = 2
$ws.Range("H133").Value = 29

# Row 172 - Islas Turcas y Caicos
$ws.Range("B172").Value = 577
$ws.Range("C172").Value = 22
$ws.Range("E172").Value = 352
$ws.Range("G172").Value = 1
$ws.Range("H172").Value = 5

# Row 183 - Mongolia
$ws.Range("B183").Value = 310
$ws.Range("C183").Value = 4
$ws.Range("E183").Value = 14

# Row 196 - Curazao
$ws.Range("B196").Value = 78
$ws.Range("C196").Value = 3
$ws.Range("E196").Value = 38
